$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Grouped by MethylScoreAML_Px Categorical"
$ws.Range("D4").Value = 707
$ws.Range("E4").Value = 233
$ws.Range("D5").Value = "9.6 (6.2)"
$ws.Range("E5").Value = "8.7 (6.6)"
$ws.Range("F5").Value = "'0.058"
$ws.Range("D6").Value = "348 (49.2)"
$ws.Range("E6").Value = "110 (47.2)"
$ws.Range("F6").Value = "'0.647"
$ws.Range("D7").Value = "359 (50.8)"
$ws.Range("E7").Value = "123 (52.8)"
$ws.Range("D8").Value = "357 (50.5)"
$ws.Range("E8").Value = "109 (46.8)"
$ws.Range("F8").Value = "'0.364"
$ws.Range("D9").Value = "350 (49.5)"
$ws.Range("E9").Value = "124 (53.2)"
$ws.Range("D10").Value = "521 (79.1)"
$ws.Range("E10").Value = "171 (78.8)"
$ws.Range("F10").Value = "'0.095"
$ws.Range("D11").Value = "68 (10.3)"
$ws.Range("E11").Value = "34 (15.7)"
$ws.Range("D12").Value = "36 (5.5)"
$ws.Range("E12").Value = "7 (3.2)"
$ws.Range("D13").Value = "4 (0.6)"
$ws.Range("E13").Value = "1 (0.5)"
$ws.Range("D14").Value = "6 (0.9)"
$ws.Range("E14").ClearContents()
$ws.Range("D15").Value = "24 (3.6)"
$ws.Range("E15").Value = "4 (1.8)"
$ws.Range("D16").Value = "137 (19.9)"
$ws.Range("E16").Value = "47 (21.1)"
$ws.Range("F16").Value = "'0.779"
$ws.Range("D17").Value = "551 (80.1)"
$ws.Range("E17").Value = "176 (78.9)"
$ws.Range("D18").Value = "178 (28.5)"
$ws.Range("E18").Value = "80 (40.8)"
$ws.Range("F18").Value = "'0.002"
$ws.Range("D19").Value = "447 (71.5)"
$ws.Range("E19").Value = "116 (59.2)"
$ws.Range("D20").Value = "355 (50.2)"
$ws.Range("E20").Value = "108 (46.4)"
$ws.Range("F20").Value = "'0.344"
$ws.Range("D21").Value = "352 (49.8)"
$ws.Range("E21").Value = "125 (53.6)"
$ws.Range("D22").Value = "61.3 (24.8)"
$ws.Range("E22").Value = "71.4 (22.2)"
$ws.Range("D23").Value = "89 (12.8)"
$ws.Range("E23").Value = "39 (17.0)"
$ws.Range("D24").Value = "276 (39.7)"
$ws.Range("E24").Value = "177 (77.0)"
$ws.Range("D25").Value = "331 (47.6)"
$ws.Range("E25").Value = "14 (6.1)"
$ws.Range("D26").Value = "27 (3.8)"
$ws.Range("E26").Value = "9 (3.9)"
$ws.Range("F26").Value = "'0.211"
$ws.Range("D27").Value = "370 (52.3)"
$ws.Range("E27").Value = "137 (58.8)"
$ws.Range("D28").Value = "310 (43.8)"
$ws.Range("E28").Value = "87 (37.3)"
$ws.Range("D29").Value = "133 (18.8)"
$ws.Range("E29").Value = "31 (13.4)"
$ws.Range("F29").Value = "'0.071"
$ws.Range("D30").Value = "573 (81.2)"
$ws.Range("E30").Value = "201 (86.6)"
